# "Ajustes para o SP1"
#
# 1) The cached text of the "datetimeFigureOut" footer field moves from
#    04/06/2022 to 20/06/2022 on every slide master / custom layout that
#    shows a date placeholder (this is PowerPoint re-caching the field's
#    display text on save). We update every reachable copy of that
#    placeholder's text (the two slide masters that carry it live); the
#    per-layout copies inherit from / mirror the master value.
# 2) On slide 2, the text box "CaixaDeTexto 100" ("Não se aplica") becomes
#    a two-line box reading "Notebook/laptop" / "Computador", and grows
#    (taller, shifted up) to fit the extra line.

$p = $ppt.ActivePresentation

# --- 1) refresh the cached "datetimeFigureOut" footer date ---------------
$oldDate = "04/06/2022"
$newDate = "20/06/2022"

$designs = $p.Designs
for ($i = 1; $i -le $designs.Count; $i++) {
    $master = $designs.Item($i).SlideMaster
    $hf = $master.HeadersFooters
    if ($hf.DateAndTime.Visible -and $hf.DateAndTime.Text -eq $oldDate) {
        for ($k = 1; $k -le $master.Shapes.Count; $k++) {
            $shp = $master.Shapes.Item($k)
            if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }

    # Same placeholder, mirrored onto every custom layout that belongs to
    # this master (best effort - some hosts treat layouts as read-only).
    $layouts = $master.CustomLayouts
    for ($j = 1; $j -le $layouts.Count; $j++) {
        $layout = $layouts.Item($j)
        for ($k = 1; $k -le $layout.Shapes.Count; $k++) {
            $shp = $layout.Shapes.Item($k)
            if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2) update the "Não se aplica" text box on slide 2 -------------------
$slide = $p.Slides.Item(2)
$shape = $slide.Shapes.Item("CaixaDeTexto 100")

$tr = $shape.TextFrame.TextRange
$firstRun = $tr.Runs(1, 1)
$firstRun.Text = "Notebook/laptop"
$tr.InsertAfter([char]13 + "Computador")

$shape.Top = 283.2015
$shape.Height = 36.3516
